$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold numeric-looking price strings that must remain text
# (matching the original t="inlineStr" cells). We force text format before
# assigning the value, then reset the style so no extra formatting is introduced.
function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextCell "D2" "89.277.09"
$ws.Range("E2").Value = "  +4.04%  "
Set-TextCell "D3" "3.296.42"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextCell "D5" "214.43"
$ws.Range("E5").Value = "  -2.51%  "
Set-TextCell "D6" "631.70"
$ws.Range("E6").Value = "  -0.49%  "
Set-TextCell "D7" "0.389"
$ws.Range("E7").Value = "  +21.01%  "
Set-TextCell "D8" "0.692"
$ws.Range("E8").Value = "  +16.53%  "
$ws.Range("E9").Value = "  -0.03%  "
Set-TextCell "D10" "3.291.60"
$ws.Range("E10").Value = "  -0.61%  "
Set-TextCell "D11" "0.581"
$ws.Range("E11").Value = "  -2.24%  "
Set-TextCell "D12" "0.187"
$ws.Range("E12").Value = "  +12.69%  "
Set-TextCell "D13" "0.0000265"
$ws.Range("E13").Value = "  -3.44%  "
Set-TextCell "D14" "34.32"
$ws.Range("E14").Value = "  +0.65%  "
Set-TextCell "D15" "3.887.35"
$ws.Range("E15").Value = "  -1.18%  "
Set-TextCell "D16" "5.39"
$ws.Range("E16").Value = "  -0.23%  "
Set-TextCell "D17" "88.922.79"
$ws.Range("E17").Value = "  +3.84%  "
Set-TextCell "D18" "3.302.79"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D19" "14.20"
$ws.Range("E19").Value = "  -2.69%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextCell "D20" "3.12"
$ws.Range("E20").Value = "  -1.43%  "
Set-TextCell "D21" "438.19"
$ws.Range("E21").Value = "  -0.62%  "
Set-TextCell "D22" "8.93"
$ws.Range("E22").Value = "  -2.21%  "
Set-TextCell "D23" "5.42"
$ws.Range("E23").Value = "  +3.45%  "
Set-TextCell "D24" "7.38"
$ws.Range("E24").Value = "  +1.14%  "
Set-TextCell "D25" "12.42"
$ws.Range("E25").Value = "  +1.53%  "
Set-TextCell "D26" "5.25"
$ws.Range("E26").Value = "  -3.52%  "
Set-TextCell "D27" "3.449.98"
Set-TextCell "D28" "77.31"
$ws.Range("E28").Value = "  -1.19%  "
Set-TextCell "D29" "0.0000135"
$ws.Range("E29").Value = "  +3.72%  "
$ws.Range("E30").Value = "  -0.02%  "
Set-TextCell "D31" "0.193"
$ws.Range("E31").Value = "  +13.70%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell "D33" "576.36"
$ws.Range("E33").Value = "  -5.58%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell "D34" "8.91"
$ws.Range("E34").Value = "  -3.56%  "
$ws.Range("E35").Value = "  -8.99%  "
Set-TextCell "D36" "7.30"
$ws.Range("E36").Value = "  +13.82%  "
Set-TextCell "D37" "1.98"
$ws.Range("E37").Value = "  -3.11%  "
Set-TextCell "D38" "0.139"
$ws.Range("E38").Value = "  -7.41%  "
Set-TextCell "D39" "22.74"
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell "D40" "21.82"
$ws.Range("E40").Value = "  +2.64%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextCell "D41" "1.00"
$ws.Range("E41").Value = "  +0.24%  "
Set-TextCell "D42" "0.401"
$ws.Range("E42").Value = "  -3.93%  "
Set-TextCell "D43" "2.04"
$ws.Range("E43").Value = "  -0.88%  "
Set-TextCell "D44" "3.01"
$ws.Range("E44").Value = "  -1.57%  "
Set-TextCell "D46" "154.31"
$ws.Range("E46").Value = "  -2.73%  "
Set-TextCell "D47" "181.35"
$ws.Range("E47").Value = "  -3.88%  "
Set-TextCell "D48" "45.15"
$ws.Range("E48").Value = "  +0.02%  "
Set-TextCell "D49" "1.31"
$ws.Range("E49").Value = "  -4.00%  "
$ws.Range("B50").Value = "Hedera"
$ws.Range("C50").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D50" "0.0680"
$ws.Range("E50").Value = "  +21.56%  "
$ws.Range("B51").Value = "Filecoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D51" "4.27"
$ws.Range("E51").Value = "  +0.31%  "
